# Apply updated odds values to the "Jogos da Semana" worksheet
# (FlashScore weekly games export) per the committed data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 10).Value = 1.03  # J2: 1.02 -> 1.03
$ws.Cells.Item(2, 11).Value = 17  # K2: 19 -> 17
$ws.Cells.Item(2, 12).Value = 1.17  # L2: 1.14 -> 1.17
$ws.Cells.Item(2, 13).Value = 5  # M2: 5.5 -> 5
$ws.Cells.Item(2, 14).Value = 1.53  # N2: 1.5 -> 1.53
$ws.Cells.Item(2, 15).Value = 2.5  # O2: 2.63 -> 2.5
$ws.Cells.Item(2, 16).Value = 1.29  # P2: 1.25 -> 1.29
$ws.Cells.Item(2, 17).Value = 3.5  # Q2: 3.75 -> 3.5
$ws.Cells.Item(2, 18).Value = 1.8  # R2: 1.75 -> 1.8
$ws.Cells.Item(2, 19).Value = 1.95  # S2: 2 -> 1.95
$ws.Cells.Item(2, 20).Value = 9.5  # T2: 10 -> 9.5
$ws.Cells.Item(2, 21).Value = 8  # U2: 8.5 -> 8
$ws.Cells.Item(2, 25).Value = 26  # Y2: 23 -> 26
$ws.Cells.Item(2, 26).Value = 17  # Z2: 19 -> 17
$ws.Cells.Item(2, 30).Value = 301  # AD2: 251 -> 301
# Row 3
$ws.Cells.Item(3, 7).Value = 2.1  # G3: 2.05 -> 2.1
$ws.Cells.Item(3, 8).Value = 3.75  # H3: 3.8 -> 3.75
$ws.Cells.Item(3, 9).Value = 3.2  # I3: 3.25 -> 3.2
$ws.Cells.Item(3, 10).Value = 1.04  # J3: 1.03 -> 1.04
$ws.Cells.Item(3, 11).Value = 13  # K3: 17 -> 13
$ws.Cells.Item(3, 12).Value = 1.22  # L3: 1.2 -> 1.22
$ws.Cells.Item(3, 13).Value = 4.33  # M3: 4.5 -> 4.33
$ws.Cells.Item(3, 14).Value = 1.67  # N3: 1.62 -> 1.67
$ws.Cells.Item(3, 15).Value = 2.2  # O3: 2.3 -> 2.2
$ws.Cells.Item(3, 16).Value = 1.33  # P3: 1.3 -> 1.33
$ws.Cells.Item(3, 17).Value = 3.25  # Q3: 3.4 -> 3.25
$ws.Cells.Item(3, 18).Value = 1.57  # R3: 1.53 -> 1.57
$ws.Cells.Item(3, 19).Value = 2.25  # S3: 2.38 -> 2.25
$ws.Cells.Item(3, 20).Value = 11  # T3: 12 -> 11
$ws.Cells.Item(3, 26).Value = 15  # Z3: 17 -> 15
$ws.Cells.Item(3, 31).Value = 13  # AE3: 15 -> 13
# Row 4
$ws.Cells.Item(4, 11).Value = 23  # K4: 26 -> 23
# Row 5
$ws.Cells.Item(5, 11).Value = 21  # K5: 19 -> 21
$ws.Cells.Item(5, 20).Value = 13  # T5: 12 -> 13
$ws.Cells.Item(5, 31).Value = 21  # AE5: 19 -> 21
# Row 6
$ws.Cells.Item(6, 8).Value = 4.2  # H6: 3.8 -> 4.2
$ws.Cells.Item(6, 9).Value = 1.55  # I6: 1.57 -> 1.55
$ws.Cells.Item(6, 22).Value = 19  # V6: 21 -> 19
$ws.Cells.Item(6, 25).Value = 41  # Y6: 51 -> 41
$ws.Cells.Item(6, 26).Value = 13  # Z6: 12 -> 13
$ws.Cells.Item(6, 27).Value = 8.5  # AA6: 8 -> 8.5
$ws.Cells.Item(6, 31).Value = 8.5  # AE6: 8 -> 8.5
$ws.Cells.Item(6, 35).Value = 13  # AI6: 15 -> 13
# Row 7
$ws.Cells.Item(7, 11).Value = 12  # K7: 13 -> 12
$ws.Cells.Item(7, 20).Value = 7.5  # T7: 7 -> 7.5
$ws.Cells.Item(7, 27).Value = 10  # AA7: 11 -> 10
$ws.Cells.Item(7, 31).Value = 19  # AE7: 21 -> 19
# Row 8
$ws.Cells.Item(8, 7).Value = 1.36  # G8: 1.38 -> 1.36
$ws.Cells.Item(8, 16).Value = 1.3  # P8: 1.33 -> 1.3
$ws.Cells.Item(8, 17).Value = 3.4  # Q8: 3.25 -> 3.4
$ws.Cells.Item(8, 20).Value = 7.5  # T8: 7 -> 7.5
$ws.Cells.Item(8, 21).Value = 7  # U8: 6.5 -> 7
$ws.Cells.Item(8, 24).Value = 11  # X8: 12 -> 11
$ws.Cells.Item(8, 26).Value = 15  # Z8: 13 -> 15
$ws.Cells.Item(8, 27).Value = 10  # AA8: 9.5 -> 10
$ws.Cells.Item(8, 29).Value = 51  # AC8: 67 -> 51
$ws.Cells.Item(8, 31).Value = 19  # AE8: 17 -> 19
# Row 10
$ws.Cells.Item(10, 7).Value = 1.38  # G10: 1.36 -> 1.38
$ws.Cells.Item(10, 8).Value = 5.25  # H10: 5 -> 5.25
$ws.Cells.Item(10, 10).Value = 1.02  # J10: 1.03 -> 1.02
$ws.Cells.Item(10, 11).Value = 19  # K10: 17 -> 19
$ws.Cells.Item(10, 16).Value = 1.25  # P10: 1.29 -> 1.25
$ws.Cells.Item(10, 17).Value = 3.75  # Q10: 3.5 -> 3.75
$ws.Cells.Item(10, 18).Value = 1.75  # R10: 1.8 -> 1.75
$ws.Cells.Item(10, 19).Value = 2  # S10: 1.95 -> 2
$ws.Cells.Item(10, 20).Value = 9  # T10: 8.5 -> 9
$ws.Cells.Item(10, 25).Value = 21  # Y10: 23 -> 21
# Row 11
$ws.Cells.Item(11, 7).Value = 3.1  # G11: 2.9 -> 3.1
$ws.Cells.Item(11, 9).Value = 2.25  # I11: 2.35 -> 2.25
$ws.Cells.Item(11, 14).Value = 2.08  # N11: 2.05 -> 2.08
$ws.Cells.Item(11, 15).Value = 1.73  # O11: 1.75 -> 1.73
$ws.Cells.Item(11, 18).Value = 1.91  # R11: 1.8 -> 1.91
$ws.Cells.Item(11, 19).Value = 1.91  # S11: 1.95 -> 1.91
$ws.Cells.Item(11, 20).Value = 9  # T11: 8.5 -> 9
$ws.Cells.Item(11, 23).Value = 34  # W11: 29 -> 34
$ws.Cells.Item(11, 24).Value = 26  # X11: 23 -> 26
$ws.Cells.Item(11, 32).Value = 10  # AF11: 11 -> 10
$ws.Cells.Item(11, 34).Value = 21  # AH11: 23 -> 21
# Row 13
$ws.Cells.Item(13, 7).Value = 2.75  # G13: 2.88 -> 2.75
$ws.Cells.Item(13, 9).Value = 2.6  # I13: 2.5 -> 2.6
$ws.Cells.Item(13, 14).Value = 1.9  # N13: 1.88 -> 1.9
$ws.Cells.Item(13, 15).Value = 1.95  # O13: 1.98 -> 1.95
$ws.Cells.Item(13, 22).Value = 10  # V13: 11 -> 10
$ws.Cells.Item(13, 34).Value = 26  # AH13: 23 -> 26
$ws.Cells.Item(13, 35).Value = 21  # AI13: 19 -> 21
# Row 14
$ws.Cells.Item(14, 7).Value = 1.27  # G14: 1.25 -> 1.27
$ws.Cells.Item(14, 8).Value = 6.5  # H14: 7 -> 6.5
$ws.Cells.Item(14, 9).Value = 8.5  # I14: 9 -> 8.5
$ws.Cells.Item(14, 11).Value = 34  # K14: 29 -> 34
$ws.Cells.Item(14, 14).Value = 1.25  # N14: 1.29 -> 1.25
$ws.Cells.Item(14, 15).Value = 4  # O14: 3.6 -> 4
$ws.Cells.Item(14, 16).Value = 1.17  # P14: 1.18 -> 1.17
$ws.Cells.Item(14, 17).Value = 5  # Q14: 4.5 -> 5
$ws.Cells.Item(14, 18).Value = 1.57  # R14: 1.62 -> 1.57
$ws.Cells.Item(14, 19).Value = 2.25  # S14: 2.2 -> 2.25
$ws.Cells.Item(14, 20).Value = 15  # T14: 13 -> 15
$ws.Cells.Item(14, 21).Value = 10  # U14: 9.5 -> 10
$ws.Cells.Item(14, 23).Value = 10  # W14: 9.5 -> 10
$ws.Cells.Item(14, 26).Value = 34  # Z14: 29 -> 34
$ws.Cells.Item(14, 30).Value = 126  # AD14: 151 -> 126
$ws.Cells.Item(14, 33).Value = 23  # AG14: 26 -> 23
# Row 16
$ws.Cells.Item(16, 7).Value = 2.5  # G16: 2.55 -> 2.5
$ws.Cells.Item(16, 9).Value = 3  # I16: 2.9 -> 3
$ws.Cells.Item(16, 10).Value = 1.11  # J16: 1.1 -> 1.11
$ws.Cells.Item(16, 11).Value = 6.5  # K16: 7 -> 6.5
$ws.Cells.Item(16, 16).Value = 1.62  # P16: 1.57 -> 1.62
$ws.Cells.Item(16, 17).Value = 2.2  # Q16: 2.25 -> 2.2
$ws.Cells.Item(16, 18).Value = 2.2  # R16: 2.1 -> 2.2
$ws.Cells.Item(16, 19).Value = 1.62  # S16: 1.67 -> 1.62
$ws.Cells.Item(16, 20).Value = 6  # T16: 6.5 -> 6
$ws.Cells.Item(16, 34).Value = 34  # AH16: 29 -> 34
# Row 17
$ws.Cells.Item(17, 10).Value = 1.07  # J17: 1.06 -> 1.07
$ws.Cells.Item(17, 11).Value = 9  # K17: 10 -> 9
$ws.Cells.Item(17, 14).Value = 2.08  # N17: 2.1 -> 2.08
$ws.Cells.Item(17, 15).Value = 1.73  # O17: 1.7 -> 1.73
$ws.Cells.Item(17, 16).Value = 1.44  # P17: 1.5 -> 1.44
$ws.Cells.Item(17, 17).Value = 2.63  # Q17: 2.5 -> 2.63
$ws.Cells.Item(17, 18).Value = 1.91  # R17: 1.95 -> 1.91
$ws.Cells.Item(17, 19).Value = 1.91  # S17: 1.8 -> 1.91
$ws.Cells.Item(17, 22).Value = 9  # V17: 9.5 -> 9
$ws.Cells.Item(17, 25).Value = 29  # Y17: 34 -> 29
$ws.Cells.Item(17, 26).Value = 9  # Z17: 8.5 -> 9
$ws.Cells.Item(17, 28).Value = 15  # AB17: 17 -> 15
$ws.Cells.Item(17, 30).Value = 301  # AD17: 351 -> 301
$ws.Cells.Item(17, 31).Value = 9.5  # AE17: 9 -> 9.5
$ws.Cells.Item(17, 32).Value = 17  # AF17: 15 -> 17
# Row 18
$ws.Cells.Item(18, 9).Value = 2.4  # I18: 2.45 -> 2.4
$ws.Cells.Item(18, 12).Value = 1.5  # L18: 1.44 -> 1.5
$ws.Cells.Item(18, 13).Value = 2.5  # M18: 2.63 -> 2.5
$ws.Cells.Item(18, 21).Value = 15  # U18: 13 -> 15
# Row 21
$ws.Cells.Item(21, 7).Value = 2.45  # G21: 2.4 -> 2.45
$ws.Cells.Item(21, 9).Value = 2.88  # I21: 2.75 -> 2.88
$ws.Cells.Item(21, 10).Value = 1.06  # J21: 1.05 -> 1.06
$ws.Cells.Item(21, 11).Value = 10  # K21: 11 -> 10
$ws.Cells.Item(21, 20).Value = 8  # T21: 8.5 -> 8
$ws.Cells.Item(21, 32).Value = 15  # AF21: 13 -> 15
$ws.Cells.Item(21, 36).Value = 34  # AJ21: 29 -> 34
# Row 22
$ws.Cells.Item(22, 7).Value = 1.83  # G22: 1.79 -> 1.83
# Row 23
$ws.Cells.Item(23, 7).Value = 2.25  # G23: 2.2 -> 2.25
$ws.Cells.Item(23, 8).Value = 3.2  # H23: 3.1 -> 3.2
$ws.Cells.Item(23, 9).Value = 3.3  # I23: 3.4 -> 3.3
$ws.Cells.Item(23, 16).Value = 1.53  # P23: 1.57 -> 1.53
$ws.Cells.Item(23, 17).Value = 2.38  # Q23: 2.25 -> 2.38
$ws.Cells.Item(23, 18).Value = 2  # R23: 2.05 -> 2
$ws.Cells.Item(23, 19).Value = 1.75  # S23: 1.7 -> 1.75
$ws.Cells.Item(23, 20).Value = 6.5  # T23: 6 -> 6.5
$ws.Cells.Item(23, 22).Value = 9.5  # V23: 10 -> 9.5
$ws.Cells.Item(23, 25).Value = 34  # Y23: 41 -> 34
$ws.Cells.Item(23, 26).Value = 7.5  # Z23: 7 -> 7.5
$ws.Cells.Item(23, 28).Value = 17  # AB23: 19 -> 17
$ws.Cells.Item(23, 30).Value = 451  # AD23: 501 -> 451
$ws.Cells.Item(23, 33).Value = 12  # AG23: 13 -> 12
$ws.Cells.Item(23, 34).Value = 34  # AH23: 41 -> 34
$ws.Cells.Item(23, 35).Value = 29  # AI23: 34 -> 29
# Row 24
$ws.Cells.Item(24, 9).Value = 2.2  # I24: 2.15 -> 2.2
# Row 25
$ws.Cells.Item(25, 9).Value = 3.75  # I25: 3.6 -> 3.75
$ws.Cells.Item(25, 11).Value = 8.5  # K25: 9 -> 8.5
$ws.Cells.Item(25, 14).Value = 2.15  # N25: 2.1 -> 2.15
$ws.Cells.Item(25, 15).Value = 1.67  # O25: 1.7 -> 1.67
$ws.Cells.Item(25, 20).Value = 6.5  # T25: 7 -> 6.5
$ws.Cells.Item(25, 21).Value = 9  # U25: 9.5 -> 9
$ws.Cells.Item(25, 23).Value = 17  # W25: 19 -> 17
$ws.Cells.Item(25, 24).Value = 17  # X25: 19 -> 17
$ws.Cells.Item(25, 28).Value = 17  # AB25: 15 -> 17
$ws.Cells.Item(25, 32).Value = 19  # AF25: 17 -> 19
# Row 26
$ws.Cells.Item(26, 18).Value = 2.5  # R26: 2.38 -> 2.5
$ws.Cells.Item(26, 19).Value = 1.5  # S26: 1.53 -> 1.5
$ws.Cells.Item(26, 24).Value = 12  # X26: 11 -> 12
$ws.Cells.Item(26, 32).Value = 67  # AF26: 81 -> 67
# Row 28
$ws.Cells.Item(28, 7).Value = 2.2  # G28: 2.1 -> 2.2
$ws.Cells.Item(28, 9).Value = 2.75  # I28: 2.9 -> 2.75
$ws.Cells.Item(28, 10).Value = 1.01  # J28: 1.02 -> 1.01
$ws.Cells.Item(28, 11).Value = 13  # K28: 12 -> 13
$ws.Cells.Item(28, 14).Value = 1.75  # N28: 1.83 -> 1.75
$ws.Cells.Item(28, 15).Value = 2.05  # O28: 2.03 -> 2.05
$ws.Cells.Item(28, 18).Value = 1.67  # R28: 1.73 -> 1.67
$ws.Cells.Item(28, 19).Value = 2.1  # S28: 2 -> 2.1
$ws.Cells.Item(28, 20).Value = 9  # T28: 8.5 -> 9
$ws.Cells.Item(28, 21).Value = 12  # U28: 11 -> 12
$ws.Cells.Item(28, 23).Value = 21  # W28: 19 -> 21
$ws.Cells.Item(28, 26).Value = 13  # Z28: 12 -> 13
$ws.Cells.Item(28, 35).Value = 21  # AI28: 23 -> 21
# Row 29
$ws.Cells.Item(29, 10).Value = 1.03  # J29: 1.04 -> 1.03
$ws.Cells.Item(29, 11).Value = 15  # K29: 13 -> 15
$ws.Cells.Item(29, 14).Value = 1.73  # N29: 1.7 -> 1.73
$ws.Cells.Item(29, 15).Value = 2.08  # O29: 2.1 -> 2.08
# Row 30
$ws.Cells.Item(30, 7).Value = 2.25  # G30: 2.35 -> 2.25
$ws.Cells.Item(30, 9).Value = 3  # I30: 2.9 -> 3
$ws.Cells.Item(30, 16).Value = 1.33  # P30: 1.36 -> 1.33
$ws.Cells.Item(30, 17).Value = 3.25  # Q30: 3 -> 3.25
$ws.Cells.Item(30, 18).Value = 1.57  # R30: 1.62 -> 1.57
$ws.Cells.Item(30, 19).Value = 2.25  # S30: 2.2 -> 2.25
$ws.Cells.Item(30, 23).Value = 21  # W30: 23 -> 21
$ws.Cells.Item(30, 24).Value = 17  # X30: 19 -> 17
$ws.Cells.Item(30, 26).Value = 13  # Z30: 12 -> 13
$ws.Cells.Item(30, 31).Value = 12  # AE30: 11 -> 12
$ws.Cells.Item(30, 32).Value = 17  # AF30: 15 -> 17
